$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.927.62"
$ws.Range("E2").Value = "  -0.88%  "

$ws.Range("D3").Value = "2.900.98"
$ws.Range("E3").Value = "  -1.12%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.00"
$ws.Range("E5").Value = "  -3.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.72"
$ws.Range("E6").Value = "  -2.38%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("E8").Value = "  -1.03%  "

$ws.Range("D9").Value = "2.899.64"
$ws.Range("E9").Value = "  -1.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.93"
$ws.Range("E10").Value = "  -3.29%  "

$ws.Range("E11").Value = "  -2.64%  "

$ws.Range("E12").Value = "  -1.31%  "

$ws.Range("E13").Value = "  -1.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.19"
$ws.Range("E14").Value = "  -0.33%  "

$ws.Range("E15").Value = "  -0.28%  "

$ws.Range("D16").Value = "3.381.33"
$ws.Range("E16").Value = "  -1.12%  "

$ws.Range("D17").Value = "61.865.08"
$ws.Range("E17").Value = "  -0.98%  "

$ws.Range("D18").Value = "2.899.04"
$ws.Range("E18").Value = "  -1.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.52"
$ws.Range("E19").Value = "  -1.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "430.38"
$ws.Range("E20").Value = "  -0.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.92"
$ws.Range("E21").Value = "  -4.11%  "

$ws.Range("E22").Value = "  -1.77%  "

$ws.Range("E23").Value = "  -1.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.99"
$ws.Range("E24").Value = "  -1.50%  "

$ws.Range("E25").Value = "  +0.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.10"
$ws.Range("E26").Value = "  -8.45%  "

$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.03"
$ws.Range("E28").Value = "  -3.13%  "

$ws.Range("E29").Value = "  +8.09%  "

$ws.Range("E30").Value = "  -3.09%  "

$ws.Range("E31").Value = "  -2.28%  "

$ws.Range("E32").Value = "  -6.01%  "

$ws.Range("E33").Value = "  -0.14%  "

$ws.Range("E34").Value = "  -2.76%  "

$ws.Range("E35").Value = "  -1.66%  "

$ws.Range("E36").Value = "  -3.28%  "

$ws.Range("E37").Value = "  -2.93%  "

$ws.Range("E38").Value = "  -1.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.85"
$ws.Range("E39").Value = "  -5.26%  "

$ws.Range("E40").Value = "  -4.37%  "

$ws.Range("E41").Value = "  -1.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.15"
$ws.Range("E42").Value = "  -2.47%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.54"
$ws.Range("E43").Value = "  +4.18%  "

$ws.Range("E44").Value = "  -2.70%  "

$ws.Range("D45").Value = "2.704.94"
$ws.Range("E45").Value = "  +0.33%  "

$ws.Range("E46").Value = "  -0.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "131.73"
$ws.Range("E47").Value = "  -2.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "345.46"
$ws.Range("E48").Value = "  -1.70%  "

$ws.Range("E50").Value = "  -0.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.57"
$ws.Range("E51").Value = "  -3.81%  "
